$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 510
$ws.Cells.Item(3, 6).Value = 1591
$ws.Cells.Item(4, 6).Value = 838
$ws.Cells.Item(5, 6).Value = 237
$ws.Cells.Item(6, 6).Value = 67
$ws.Cells.Item(7, 6).Value = 1146
$ws.Cells.Item(8, 6).Value = 749
$ws.Cells.Item(9, 6).Value = 794
$ws.Cells.Item(10, 6).Value = 1442
$ws.Cells.Item(11, 6).Value = 289
$ws.Cells.Item(12, 6).Value = 1034
$ws.Cells.Item(13, 6).Value = 32
$ws.Cells.Item(14, 6).Value = 66
$ws.Cells.Item(15, 6).Value = 193
$ws.Cells.Item(16, 6).Value = 50
$ws.Cells.Item(17, 6).Value = 473
$ws.Cells.Item(18, 6).Value = 23
$ws.Cells.Item(19, 6).Value = 24
$ws.Cells.Item(20, 6).Value = 3
$ws.Cells.Item(22, 6).Value = 297
$ws.Cells.Item(23, 6).Value = 551
$ws.Cells.Item(24, 6).Value = 568
$ws.Cells.Item(25, 6).Value = 757
$ws.Cells.Item(26, 6).Value = 245
$ws.Cells.Item(27, 6).Value = 177
$ws.Cells.Item(28, 6).Value = 370

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 9
$ws.Cells.Item(3, 6).Value = 1000
$ws.Cells.Item(5, 6).Value = 268
$ws.Cells.Item(7, 6).Value = 145
$ws.Cells.Item(8, 6).Value = 66
$ws.Cells.Item(9, 6).Value = 586
$ws.Cells.Item(10, 6).Value = 80

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 234

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 510
$ws.Cells.Item(3, 6).Value = 234
$ws.Cells.Item(4, 6).Value = 1591
$ws.Cells.Item(5, 6).Value = 9
$ws.Cells.Item(6, 6).Value = 838
$ws.Cells.Item(7, 6).Value = 237
$ws.Cells.Item(8, 6).Value = 1000
$ws.Cells.Item(9, 6).Value = 67
$ws.Cells.Item(10, 6).Value = 1146
$ws.Cells.Item(11, 6).Value = 749
$ws.Cells.Item(12, 6).Value = 794
$ws.Cells.Item(13, 6).Value = 1442
$ws.Cells.Item(14, 6).Value = 289
$ws.Cells.Item(15, 6).Value = 1034
$ws.Cells.Item(16, 6).Value = 32
$ws.Cells.Item(17, 6).Value = 66
$ws.Cells.Item(18, 6).Value = 193
$ws.Cells.Item(19, 6).Value = 50
$ws.Cells.Item(20, 6).Value = 473
$ws.Cells.Item(21, 6).Value = 23
$ws.Cells.Item(22, 6).Value = 24
$ws.Cells.Item(24, 6).Value = 3
$ws.Cells.Item(25, 6).Value = 268
$ws.Cells.Item(27, 6).Value = 297
$ws.Cells.Item(29, 6).Value = 145
$ws.Cells.Item(30, 6).Value = 145
$ws.Cells.Item(31, 6).Value = 551
$ws.Cells.Item(32, 6).Value = 568
$ws.Cells.Item(33, 6).Value = 757
$ws.Cells.Item(34, 6).Value = 245
$ws.Cells.Item(35, 6).Value = 66
$ws.Cells.Item(36, 6).Value = 177
$ws.Cells.Item(37, 6).Value = 586
$ws.Cells.Item(38, 6).Value = 80
$ws.Cells.Item(39, 6).Value = 80
$ws.Cells.Item(41, 6).Value = 370
